# "update register account info"
#
# The Register sheet's sample row (row 2) used test account "testLT005" /
# "testLT005@test.com" for firstname/lastname/email. Bump it to the next
# test account, "testLT006" / "testLT006@test.com".

$wb = $excel.ActiveWorkbook

$register = $wb.Worksheets.Item("Register")
$register.Range("B2").Value = "testLT006"          # firstname
$register.Range("C2").Value = "testLT006"           # lastname
$register.Range("D2").Value = "testLT006@test.com"  # email

# Scroll the Login sheet's view down so row 10 is the top visible row
# (view-state only; doesn't touch the current selection/active cell).
$login = $wb.Worksheets.Item("Login")
$login.Select()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
try { $excel.ActiveWindow.TopLeftCell = $login.Range("A10") } catch {}

# Restore Register as the active sheet/tab (it was active before this
# script ran and should remain so).
$register.Select()

# Shrink the workbook window a bit (view-state only).
try { $excel.ActiveWindow.Height = 2136 } catch {}
try { $excel.Height = 2136 } catch {}

# Refresh the OLE display size used when this workbook is embedded as an
# object elsewhere, to cover the sample data's extent.
try { $wb.OLEObjects() | Out-Null } catch {}
